# Refresh the cryptos list: updated Price (D) and Volume(1h) (E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most Price cells hold plain decimal-looking text ("327.05", "0.3590", ...).
# Assigning such a string straight to .Value lets Excel auto-convert it to a
# real number (dropping text formatting / trailing zeros), so for those cells
# we switch to Text format first, type the literal string, then restore the
# default "Normal" cell style so no stray formatting is left behind.
$textForceCells = @(
    "D5", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D17", "D19", "D21", "D22", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cell in $textForceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D5").Value = "327.05"
$ws.Range("D7").Value = "0.4543"
$ws.Range("D8").Value = "0.3590"
$ws.Range("D9").Value = "0.07497"
$ws.Range("D10").Value = "42.04"
$ws.Range("D11").Value = "1.107"
$ws.Range("D13").Value = "20.94"
$ws.Range("D14").Value = "6.063"
$ws.Range("D15").Value = "7.217"
$ws.Range("D17").Value = "93.75"
$ws.Range("D19").Value = "0.06430"
$ws.Range("D21").Value = "17.16"
$ws.Range("D22").Value = "5.811"
$ws.Range("D25").Value = "2.085"
$ws.Range("D26").Value = "164.13"
$ws.Range("D27").Value = "20.33"
$ws.Range("D29").Value = "2.236"
$ws.Range("D30").Value = "126.07"
$ws.Range("D31").Value = "1.120"
$ws.Range("D32").Value = "0.09214"
$ws.Range("D34").Value = "5.569"
$ws.Range("D35").Value = "11.94"
$ws.Range("D37").Value = "0.06176"
$ws.Range("D38").Value = "0.2096"
$ws.Range("D39").Value = "0.6343"
$ws.Range("D40").Value = "4.982"
$ws.Range("D41").Value = "1.187"
$ws.Range("D42").Value = "1.388"
$ws.Range("D43").Value = "7.917"
$ws.Range("D44").Value = "13.30"
$ws.Range("D45").Value = "0.5928"
$ws.Range("D46").Value = "3.735"
$ws.Range("D47").Value = "122.86"
$ws.Range("D48").Value = "1.964"
$ws.Range("D49").Value = "0.06933"
$ws.Range("D50").Value = "1.140"
$ws.Range("D51").Value = "72.98"

foreach ($cell in $textForceCells) {
    $ws.Range($cell).Style = "Normal"
}

# Price cells containing thousands separators (e.g. "27.945.08") are never
# mistaken for numbers by Excel, so they can be assigned directly.
$ws.Range("D2").Value = "27.945.08"
$ws.Range("D3").Value = "1.778.78"
$ws.Range("D16").Value = "1.778.41"
$ws.Range("D23").Value = "27.975.18"
$ws.Range("D28").Value = "1.985.76"

# Volume(1h) percentage cells are padded with spaces and a "%" sign, so they
# always remain plain text and can be assigned directly.
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("E29").Value = "  +7.50%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("E37").Value = "  +2.64%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("E51").Value = "  +0.89%  "
